$d = $word.ActiveDocument

$d.Content.Find.Execute("33-17=", $true, $false, $false, $false, $false, $true, 1, $false, "87-31=", 2) | Out-Null
$d.Content.Find.Execute("72-63=", $true, $false, $false, $false, $false, $true, 1, $false, "86-42=", 2) | Out-Null
$d.Content.Find.Execute("36+62=", $true, $false, $false, $false, $false, $true, 1, $false, "63-14=", 2) | Out-Null
$d.Content.Find.Execute("31+32=", $true, $false, $false, $false, $false, $true, 1, $false, "66-29=", 2) | Out-Null
$d.Content.Find.Execute("91-29=", $true, $false, $false, $false, $false, $true, 1, $false, "5+22=", 2) | Out-Null
$d.Content.Find.Execute("85-12=", $true, $false, $false, $false, $false, $true, 1, $false, "26-22=", 2) | Out-Null
$d.Content.Find.Execute("21-6=", $true, $false, $false, $false, $false, $true, 1, $false, "99-26=", 2) | Out-Null
$d.Content.Find.Execute("51+41=", $true, $false, $false, $false, $false, $true, 1, $false, "89-47=", 2) | Out-Null
$d.Content.Find.Execute("67-5=", $true, $false, $false, $false, $false, $true, 1, $false, "54-48=", 2) | Out-Null
$d.Content.Find.Execute("34+14=", $true, $false, $false, $false, $false, $true, 1, $false, "34+13=", 2) | Out-Null
$d.Content.Find.Execute("15+55=", $true, $false, $false, $false, $false, $true, 1, $false, "42-26=", 2) | Out-Null
$d.Content.Find.Execute("3+51=", $true, $false, $false, $false, $false, $true, 1, $false, "39+21=", 2) | Out-Null
$d.Content.Find.Execute("74-22=", $true, $false, $false, $false, $false, $true, 1, $false, "91-51=", 2) | Out-Null
$d.Content.Find.Execute("80-0=", $true, $false, $false, $false, $false, $true, 1, $false, "15-8=", 2) | Out-Null
$d.Content.Find.Execute("33+6=", $true, $false, $false, $false, $false, $true, 1, $false, "47-20=", 2) | Out-Null
$d.Content.Find.Execute("36+30=", $true, $false, $false, $false, $false, $true, 1, $false, "42-10=", 2) | Out-Null
$d.Content.Find.Execute("56-52=", $true, $false, $false, $false, $false, $true, 1, $false, "43-33=", 2) | Out-Null
$d.Content.Find.Execute("28+13=", $true, $false, $false, $false, $false, $true, 1, $false, "11+48=", 2) | Out-Null
$d.Content.Find.Execute("89-60=", $true, $false, $false, $false, $false, $true, 1, $false, "79-36=", 2) | Out-Null
$d.Content.Find.Execute("83-46=", $true, $false, $false, $false, $false, $true, 1, $false, "88-7=", 2) | Out-Null
$d.Content.Find.Execute("67-21=", $true, $false, $false, $false, $false, $true, 1, $false, "14+51=", 2) | Out-Null
$d.Content.Find.Execute("1-0=", $true, $false, $false, $false, $false, $true, 1, $false, "55-1=", 2) | Out-Null
$d.Content.Find.Execute("81-36=", $true, $false, $false, $false, $false, $true, 1, $false, "44+17=", 2) | Out-Null
$d.Content.Find.Execute("9+78=", $true, $false, $false, $false, $false, $true, 1, $false, "33+59=", 2) | Out-Null
$d.Content.Find.Execute("22+18=", $true, $false, $false, $false, $false, $true, 1, $false, "39-39=", 2) | Out-Null
$d.Content.Find.Execute("55-38=", $true, $false, $false, $false, $false, $true, 1, $false, "23+66=", 2) | Out-Null
$d.Content.Find.Execute("36+59=", $true, $false, $false, $false, $false, $true, 1, $false, "8+52=", 2) | Out-Null
$d.Content.Find.Execute("54+19=", $true, $false, $false, $false, $false, $true, 1, $false, "0+90=", 2) | Out-Null
$d.Content.Find.Execute("25+53=", $true, $false, $false, $false, $false, $true, 1, $false, "5+85=", 2) | Out-Null
$d.Content.Find.Execute("76-23=", $true, $false, $false, $false, $false, $true, 1, $false, "85+14=", 2) | Out-Null
$d.Content.Find.Execute("73-6=", $true, $false, $false, $false, $false, $true, 1, $false, "47-0=", 2) | Out-Null
$d.Content.Find.Execute("56-5=", $true, $false, $false, $false, $false, $true, 1, $false, "41+17=", 2) | Out-Null
$d.Content.Find.Execute("74-5=", $true, $false, $false, $false, $false, $true, 1, $false, "17+2=", 2) | Out-Null
$d.Content.Find.Execute("11+49=", $true, $false, $false, $false, $false, $true, 1, $false, "44+27=", 2) | Out-Null
$d.Content.Find.Execute("48+40=", $true, $false, $false, $false, $false, $true, 1, $false, "12+70=", 2) | Out-Null
$d.Content.Find.Execute("64+0=", $true, $false, $false, $false, $false, $true, 1, $false, "40-35=", 2) | Out-Null
$d.Content.Find.Execute("35+48=", $true, $false, $false, $false, $false, $true, 1, $false, "69-40=", 2) | Out-Null
$d.Content.Find.Execute("95-92=", $true, $false, $false, $false, $false, $true, 1, $false, "40+29=", 2) | Out-Null
$d.Content.Find.Execute("53-25=", $true, $false, $false, $false, $false, $true, 1, $false, "21+9=", 2) | Out-Null
$d.Content.Find.Execute("71-54=", $true, $false, $false, $false, $false, $true, 1, $false, "9+84=", 2) | Out-Null
$d.Content.Find.Execute("93-45=", $true, $false, $false, $false, $false, $true, 1, $false, "96-73=", 2) | Out-Null
$d.Content.Find.Execute("94-52=", $true, $false, $false, $false, $false, $true, 1, $false, "56-50=", 2) | Out-Null
$d.Content.Find.Execute("35+55=", $true, $false, $false, $false, $false, $true, 1, $false, "6+47=", 2) | Out-Null
$d.Content.Find.Execute("15+25=", $true, $false, $false, $false, $false, $true, 1, $false, "99-60=", 2) | Out-Null
$d.Content.Find.Execute("80+17=", $true, $false, $false, $false, $false, $true, 1, $false, "17+30=", 2) | Out-Null
$d.Content.Find.Execute("91-28=", $true, $false, $false, $false, $false, $true, 1, $false, "53-31=", 2) | Out-Null
$d.Content.Find.Execute("94+4=", $true, $false, $false, $false, $false, $true, 1, $false, "37+13=", 2) | Out-Null
$d.Content.Find.Execute("35+59=", $true, $false, $false, $false, $false, $true, 1, $false, "70+26=", 2) | Out-Null
$d.Content.Find.Execute("41-34=", $true, $false, $false, $false, $false, $true, 1, $false, "50+43=", 2) | Out-Null
$d.Content.Find.Execute("52+3=", $true, $false, $false, $false, $false, $true, 1, $false, "85-60=", 2) | Out-Null
$d.Content.Find.Execute("42+18=", $true, $false, $false, $false, $false, $true, 1, $false, "7+25=", 2) | Out-Null
$d.Content.Find.Execute("14+3=", $true, $false, $false, $false, $false, $true, 1, $false, "67-36=", 2) | Out-Null
$d.Content.Find.Execute("35+20=", $true, $false, $false, $false, $false, $true, 1, $false, "53-31=", 2) | Out-Null
$d.Content.Find.Execute("30-3=", $true, $false, $false, $false, $false, $true, 1, $false, "38-5=", 2) | Out-Null
$d.Content.Find.Execute("64+18=", $true, $false, $false, $false, $false, $true, 1, $false, "15+43=", 2) | Out-Null
$d.Content.Find.Execute("59-2=", $true, $false, $false, $false, $false, $true, 1, $false, "39+40=", 2) | Out-Null
$d.Content.Find.Execute("49+40=", $true, $false, $false, $false, $false, $true, 1, $false, "58-39=", 2) | Out-Null
$d.Content.Find.Execute("50-18=", $true, $false, $false, $false, $false, $true, 1, $false, "3+42=", 2) | Out-Null
$d.Content.Find.Execute("23-2=", $true, $false, $false, $false, $false, $true, 1, $false, "24+18=", 2) | Out-Null
$d.Content.Find.Execute("53-30=", $true, $false, $false, $false, $false, $true, 1, $false, "93-63=", 2) | Out-Null
$d.Content.Find.Execute("71-37=", $true, $false, $false, $false, $false, $true, 1, $false, "47+13=", 2) | Out-Null
$d.Content.Find.Execute("48+13=", $true, $false, $false, $false, $false, $true, 1, $false, "85-51=", 2) | Out-Null
$d.Content.Find.Execute("61-57=", $true, $false, $false, $false, $false, $true, 1, $false, "99-69=", 2) | Out-Null
$d.Content.Find.Execute("71-26=", $true, $false, $false, $false, $false, $true, 1, $false, "67+13=", 2) | Out-Null
$d.Content.Find.Execute("28+57=", $true, $false, $false, $false, $false, $true, 1, $false, "5+10=", 2) | Out-Null
$d.Content.Find.Execute("31+35=", $true, $false, $false, $false, $false, $true, 1, $false, "41-36=", 2) | Out-Null
$d.Content.Find.Execute("4+68=", $true, $false, $false, $false, $false, $true, 1, $false, "3+20=", 2) | Out-Null
$d.Content.Find.Execute("2+22=", $true, $false, $false, $false, $false, $true, 1, $false, "31+62=", 2) | Out-Null
$d.Content.Find.Execute("86-37=", $true, $false, $false, $false, $false, $true, 1, $false, "44-5=", 2) | Out-Null
$d.Content.Find.Execute("94-91=", $true, $false, $false, $false, $false, $true, 1, $false, "27-18=", 2) | Out-Null
$d.Content.Find.Execute("86-84=", $true, $false, $false, $false, $false, $true, 1, $false, "11+55=", 2) | Out-Null
$d.Content.Find.Execute("83-49=", $true, $false, $false, $false, $false, $true, 1, $false, "24-7=", 2) | Out-Null
$d.Content.Find.Execute("93-31=", $true, $false, $false, $false, $false, $true, 1, $false, "84-35=", 2) | Out-Null
$d.Content.Find.Execute("49-18=", $true, $false, $false, $false, $false, $true, 1, $false, "77-15=", 2) | Out-Null
$d.Content.Find.Execute("55-26=", $true, $false, $false, $false, $false, $true, 1, $false, "21+5=", 2) | Out-Null
$d.Content.Find.Execute("90-17=", $true, $false, $false, $false, $false, $true, 1, $false, "42-16=", 2) | Out-Null
$d.Content.Find.Execute("50-19=", $true, $false, $false, $false, $false, $true, 1, $false, "89-42=", 2) | Out-Null
$d.Content.Find.Execute("28-27=", $true, $false, $false, $false, $false, $true, 1, $false, "64-44=", 2) | Out-Null
$d.Content.Find.Execute("22-19=", $true, $false, $false, $false, $false, $true, 1, $false, "52+47=", 2) | Out-Null
$d.Content.Find.Execute("40+36=", $true, $false, $false, $false, $false, $true, 1, $false, "46-3=", 2) | Out-Null
$d.Content.Find.Execute("13-6=", $true, $false, $false, $false, $false, $true, 1, $false, "30+57=", 2) | Out-Null
$d.Content.Find.Execute("90-59=", $true, $false, $false, $false, $false, $true, 1, $false, "75+5=", 2) | Out-Null
$d.Content.Find.Execute("67-34=", $true, $false, $false, $false, $false, $true, 1, $false, "64+7=", 2) | Out-Null
$d.Content.Find.Execute("93+3=", $true, $false, $false, $false, $false, $true, 1, $false, "44-0=", 2) | Out-Null
$d.Content.Find.Execute("20+32=", $true, $false, $false, $false, $false, $true, 1, $false, "66-30=", 2) | Out-Null
$d.Content.Find.Execute("78-65=", $true, $false, $false, $false, $false, $true, 1, $false, "69+24=", 2) | Out-Null
$d.Content.Find.Execute("75-52=", $true, $false, $false, $false, $false, $true, 1, $false, "38+5=", 2) | Out-Null
$d.Content.Find.Execute("5+73=", $true, $false, $false, $false, $false, $true, 1, $false, "62+18=", 2) | Out-Null
$d.Content.Find.Execute("94-88=", $true, $false, $false, $false, $false, $true, 1, $false, "21+29=", 2) | Out-Null
$d.Content.Find.Execute("48+24=", $true, $false, $false, $false, $false, $true, 1, $false, "6+64=", 2) | Out-Null
$d.Content.Find.Execute("82-47=", $true, $false, $false, $false, $false, $true, 1, $false, "45+44=", 2) | Out-Null
$d.Content.Find.Execute("37+51=", $true, $false, $false, $false, $false, $true, 1, $false, "49-31=", 2) | Out-Null
$d.Content.Find.Execute("84+13=", $true, $false, $false, $false, $false, $true, 1, $false, "64-28=", 2) | Out-Null
$d.Content.Find.Execute("33-7=", $true, $false, $false, $false, $false, $true, 1, $false, "94-55=", 2) | Out-Null
$d.Content.Find.Execute("62+32=", $true, $false, $false, $false, $false, $true, 1, $false, "30+15=", 2) | Out-Null
$d.Content.Find.Execute("12+29=", $true, $false, $false, $false, $false, $true, 1, $false, "9+50=", 2) | Out-Null
$d.Content.Find.Execute("76-12=", $true, $false, $false, $false, $false, $true, 1, $false, "66-29=", 2) | Out-Null
$d.Content.Find.Execute("53+26=", $true, $false, $false, $false, $false, $true, 1, $false, "34+25=", 2) | Out-Null
$d.Content.Find.Execute("13-8=", $true, $false, $false, $false, $false, $true, 1, $false, "14+72=", 2) | Out-Null
$d.Content.Find.Execute("67+17=", $true, $false, $false, $false, $false, $true, 1, $false, "36+7=", 2) | Out-Null
